# fix(publipostage): Add space before ":"
#
# The "statut_name" column (B) held labels like "4: ..." — a colon glued
# directly to the digit. Add a space before the colon for every status
# label. Also two "intervention_type" entries (I11/I12), previously
# BEHAVIORAL, are corrected to DIAGNOSTIC_TEST.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "4: pas de résultats postés ni publiés"              = "4 : pas de résultats postés ni publiés"
    "2: résultats postés ou publiés entre 12 et 36 mois" = "2 : résultats postés ou publiés entre 12 et 36 mois"
    "3: résultats postés ou publiés après les 36 mois"   = "3 : résultats postés ou publiés après les 36 mois"
}

# statut_name lives in column B, data rows 2-14
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null -and $replacements.ContainsKey($val)) {
        $cell.Value = $replacements[$val]
    }
}

# intervention_type (column I) fix: rows 11-12 go from BEHAVIORAL to DIAGNOSTIC_TEST
$ws.Range("I11").Value = "DIAGNOSTIC_TEST"
$ws.Range("I12").Value = "DIAGNOSTIC_TEST"
